$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data (and BabyDogeCoin/RocketPoolETH row swap)

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.998.75"
$ws.Range("E2").Value = "  +0.13%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.641.43"
$ws.Range("E3").Value = "  -0.45%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.75%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.20"
$ws.Range("E5").Value = "  -0.39%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5078"
$ws.Range("E6").Value = "  -0.60%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.002"
$ws.Range("E7").Value = "  -0.57%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2583"
$ws.Range("E8").Value = "  -0.06%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06366"
$ws.Range("E9").Value = "  -0.96%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.86"
$ws.Range("E10").Value = "  +0.90%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07764"
$ws.Range("E11").Value = "  -0.33%  "

# Row 12
$ws.Range("E12").Value = "  -1.06%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.629.12"
$ws.Range("E13").Value = "  -1.37%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5486"
$ws.Range("E14").Value = "  +0.23%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0₅7754"
$ws.Range("E15").Value = "  -1.83%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.40"
$ws.Range("E16").Value = "  -0.69%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.010.47"
$ws.Range("E17").Value = "  -0.09%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.001"
$ws.Range("E18").Value = "  -0.67%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "196.70"
$ws.Range("E19").Value = "  -1.02%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.443"
$ws.Range("E20").Value = "  -0.60%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.943"
$ws.Range("E21").Value = "  -0.88%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.104"
$ws.Range("E22").Value = "  +0.40%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.003"
$ws.Range("E23").Value = "  -0.56%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.895"

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.11"
$ws.Range("E25").Value = "  +2.71%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1239"
$ws.Range("E26").Value = "  +7.61%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.880"
$ws.Range("E27").Value = "  -0.52%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.62"
$ws.Range("E28").Value = "  -1.01%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.239"
$ws.Range("E29").Value = "  -0.43%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.04886"
$ws.Range("E30").Value = "  -2.88%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.279"
$ws.Range("E31").Value = "  -0.52%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.227"
$ws.Range("E32").Value = "  +0.58%  "

# Row 33
$ws.Range("E33").Value = "  -0.13%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.375"
$ws.Range("E34").Value = "  +0.39%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9161"
$ws.Range("E35").Value = "  +2.28%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.573"
$ws.Range("E36").Value = "  -0.96%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5551"
$ws.Range("E37").Value = "  +0.08%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.088.53"
$ws.Range("E38").Value = "  -4.43%  "

# Row 39
$ws.Range("E39").Value = "  +0.36%  "

# Row 40
$ws.Range("E40").Value = "  -0.66%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.524"
$ws.Range("E41").Value = "  -1.59%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.606"
$ws.Range("E42").Value = "  -1.03%  "

# Row 43
$ws.Range("E43").Value = "  -1.60%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.25"
$ws.Range("E44").Value = "  -0.88%  "

# Row 45
$ws.Range("B45").Value = "BabyDogeCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0₈120"
$ws.Range("E45").Value = "  -4.35%  "

# Row 46
$ws.Range("B46").Value = "RocketPoolETH"
$ws.Range("C46").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.779.72"
$ws.Range("E46").Value = "  -0.32%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4539"
$ws.Range("E47").Value = "  -0.13%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "55.53"
$ws.Range("E48").Value = "  +0.37%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.004"
$ws.Range("E49").Value = "  -0.35%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05226"
$ws.Range("E50").Value = "  +2.54%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.567"
$ws.Range("E51").Value = "  +2.06%  "
